$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "September 19, 2025", $true, $false, $false, $false, $false,
    $true, 1, $false, "September 21, 2025", 2) | Out-Null

# -----------------------------------------------------------------
# 2. Split the mailing address paragraph into two lines and add a
#    blank line after it:
#    "20635 Maria Court, Castro Valley CA 94546"
#    becomes
#    "20635 Maria Court"
#    "Castro Valley, CA 94546"
#    <blank paragraph>
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "20635 Maria Court, Castro Valley CA 94546", $true, $false, $false, $false, $false,
    $true, 1, $false, "20635 Maria Court^pCastro Valley, CA 94546^p", 2) | Out-Null

# -----------------------------------------------------------------
# 3. Remove the two blank paragraphs that immediately follow the
#    "Board of Directors" line (a No Spacing blank paragraph and a
#    Title-styled blank paragraph), leaving the next blank Title
#    paragraph in place.
# -----------------------------------------------------------------
$boardPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Board of Directors*") {
        $boardPara = $p
        break
    }
}

if ($boardPara -ne $null) {
    $markPos = $boardPara.Range.End
    $d.Range($markPos, $markPos + 1).Delete() | Out-Null
    $d.Range($markPos, $markPos + 1).Delete() | Out-Null
}
